# Applies the "Unique Data" worksheet addition + header highlight to
# APT_NGINMessage.xlsx, matching the supplied OOXML diff.

$wb = $excel.ActiveWorkbook
$wsMain = $wb.ActiveSheet

# Highlight the "NGINMessage_SANNumber" header cell (C1) in yellow to mark
# it as the column that must contain unique data.
$headerCell = $wsMain.Range("C1")
$headerCell.Interior.Color = 65535

# Move/refresh the active selection on the main sheet (C5 -> C3 per diff).
$wsMain.Range("C3").Select() | Out-Null

# Add the new "Unique Data" worksheet after the existing sheet.
$wsUnique = $wb.Worksheets.Add([System.Type]::Missing, $wsMain)
$wsUnique.Name = "Unique Data"

# Populate the new sheet.
$wsUnique.Range("A1").Value = "Columns for which unique data to be provided"
$wsUnique.Range("A2").Value = "NGINMessage_SANNumber"

# Style the header cell with the same fill used for the other section
# headers in the workbook (theme-8 fill, ~40% tint), matching cellXf
# index 6 in the diff -- only the fill is applied (default font/format).
# Read the resolved color straight off an existing themed header cell so
# the new fill matches exactly, regardless of theme palette specifics.
$wsUnique.Range("A1").Interior.Color = $wsMain.Range("A1").Interior.Color

# Column width / best fit to match <col min="1" max="1" width="43" bestFit="1".../>
$wsUnique.Columns("A").ColumnWidth = 42.166666666666664

$wsUnique.Range("A2").Select() | Out-Null

$wsMain.Activate() | Out-Null
